$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "edit1"
$ws.Range("B6").Value = "riya-morankar"
$ws.Range("C6").Value = "Merged"

# E6 holds a date-shaped string ("2025-06-18") that must stay plain text,
# matching how the other rows in this column store their dates (as text,
# not as a numeric date serial). Force text interpretation via NumberFormat,
# then restore the cell to the default "Normal" style so no stray
# number-format style lingers on the cell.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2025-06-18"
$ws.Range("E6").Style = "Normal"

$ws.Range("F6").Value = "N/A"
